$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Andesite / Blackstone / Diorite / Tuff" block up by one,
# replacing the old "Mud" entry at B12 (the recipe list is being reordered and
# a new "Dragonrot Dipping" recipe section appended below).
$ws.Range("B12").Value = "Andesite"
$ws.Range("B13").Value = "Diorite"
$ws.Range("B14").Value = "Tuff"

# Remove the now-superseded rows that used to hold this tail of the list.
$ws.Range("B15:B17").ClearContents()

# New recipe section: "Dragonrot Dipping" consuming Soul Sand, Blackstone and Mud.
$ws.Range("B22").Value = "Dragonrot Dipping:"
$ws.Range("B23").Value = "Soul Sand"
$ws.Range("B24").Value = "Blackstone"
$ws.Range("B25").Value = "Mud"

# Match the saved selection state recorded in the workbook.
$ws.Range("B12:B14").Select()
